# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (shape 2) gets a new table style GUID.
# 2) The presentation's design theme (ppt/theme/theme2.xml, the theme
#    actually applied to the slide master / all slides) is swapped from
#    the "Integral" / "Red Violet" palette back to the stock
#    "Office Theme" / "Office" palette. PowerPoint's Theme.Name /
#    ColorScheme.Name are read-only in the object model, so we recreate
#    the swap by writing each of the twelve theme colours individually
#    (fonts/format scheme are already identical between the two themes).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 --------------------------------------
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{6A886D21-1712-4B03-8B1D-E20337219FB6}")

# --- 2) Swap the live theme's colour scheme back to "Office" --------
# RGB is a COLORREF (0x00BBGGRR), so pass R + G*256 + B*65536 for a
# target RRGGBB hex colour - this mirrors VBA's RGB(r,g,b) helper.
function HexColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Colors(1).RGB  = HexColor 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = HexColor 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = HexColor 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = HexColor 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = HexColor 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = HexColor 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = HexColor 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = HexColor 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = HexColor 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = HexColor 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = HexColor 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = HexColor 0x95 0x4F 0x72   # folHlink
